$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (hyp -> Sheet1)
$ws.Name = "Sheet1"

# Header row (row 1), columns B..L - reordered labels
$ws.Range("B1").Value2 = "Convolutional"
$ws.Range("C1").Value2 = "DNN"
$ws.Range("D1").Value2 = "Feedforward"
$ws.Range("E1").Value2 = "Sequence"
$ws.Range("F1").Value2 = "Attention"
$ws.Range("G1").Value2 = "Embedding"
$ws.Range("H1").Value2 = "Other"
$ws.Range("I1").Value2 = "word2vec"
$ws.Range("J1").Value2 = "Graph"
$ws.Range("K1").Value2 = "Deep Belief Network"
$ws.Range("L1").Value2 = "total"

# Data rows (category rows reordered/renamed, values updated, new total column L)
$ws.Range("A2").Value2 = "Clone Detection"
$ws.Range("B2").Value2 = 2
$ws.Range("C2").Value2 = 1
$ws.Range("D2").Value2 = 1
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0
$ws.Range("G2").Value2 = 0
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 0
$ws.Range("L2").Value2 = 6

$ws.Range("A3").Value2 = "Code Synthesis"
$ws.Range("B3").Value2 = 4
$ws.Range("C3").Value2 = 0
$ws.Range("D3").Value2 = 0
$ws.Range("E3").Value2 = 16
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1
$ws.Range("H3").Value2 = 1
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("L3").Value2 = 23

$ws.Range("A4").Value2 = "Feature Envy Detection"
$ws.Range("B4").Value2 = 1
$ws.Range("C4").Value2 = 0
$ws.Range("D4").Value2 = 0
$ws.Range("E4").Value2 = 0
$ws.Range("F4").Value2 = 0
$ws.Range("G4").Value2 = 0
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 1

$ws.Range("A5").Value2 = "Program Repair"
$ws.Range("B5").Value2 = 0
$ws.Range("C5").Value2 = 0
$ws.Range("D5").Value2 = 0
$ws.Range("E5").Value2 = 5
$ws.Range("F5").Value2 = 0
$ws.Range("G5").Value2 = 0
$ws.Range("H5").Value2 = 1
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 6

$ws.Range("A6").Value2 = "Software Categorization"
$ws.Range("B6").Value2 = 0
$ws.Range("C6").Value2 = 1
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 0
$ws.Range("L6").Value2 = 2

$ws.Range("A7").Value2 = "Software Energy Metrics"
$ws.Range("B7").Value2 = 0
$ws.Range("C7").Value2 = 0
$ws.Range("D7").Value2 = 1
$ws.Range("E7").Value2 = 0
$ws.Range("F7").Value2 = 0
$ws.Range("G7").Value2 = 0
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 0
$ws.Range("K7").Value2 = 0
$ws.Range("L7").Value2 = 1

$ws.Range("A8").Value2 = "Testing"
$ws.Range("B8").Value2 = 0
$ws.Range("C8").Value2 = 0
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = 0
$ws.Range("F8").Value2 = 0
$ws.Range("G8").Value2 = 0
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = 1
$ws.Range("J8").Value2 = 0
$ws.Range("K8").Value2 = 0
$ws.Range("L8").Value2 = 1

$ws.Range("A9").Value2 = "Vulnerability Detection"
$ws.Range("B9").Value2 = 1
$ws.Range("C9").Value2 = 0
$ws.Range("D9").Value2 = 0
$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = 0
$ws.Range("G9").Value2 = 0
$ws.Range("H9").Value2 = 0
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 1
$ws.Range("K9").Value2 = 0
$ws.Range("L9").Value2 = 4

$ws.Range("A10").Value2 = "bug localization"
$ws.Range("B10").Value2 = 1
$ws.Range("C10").Value2 = 0
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 0
$ws.Range("F10").Value2 = 0
$ws.Range("G10").Value2 = 0
$ws.Range("H10").Value2 = 1
$ws.Range("I10").Value2 = 0
$ws.Range("J10").Value2 = 0
$ws.Range("K10").Value2 = 0
$ws.Range("L10").Value2 = 2

$ws.Range("A11").Value2 = "code comprehension"
$ws.Range("B11").Value2 = 2
$ws.Range("C11").Value2 = 0
$ws.Range("D11").Value2 = 0
$ws.Range("E11").Value2 = 11
$ws.Range("F11").Value2 = 0
$ws.Range("G11").Value2 = 1
$ws.Range("H11").Value2 = 0
$ws.Range("I11").Value2 = 1
$ws.Range("J11").Value2 = 1
$ws.Range("K11").Value2 = 0
$ws.Range("L11").Value2 = 16

$ws.Range("A12").Value2 = "code smell"
$ws.Range("B12").Value2 = 0
$ws.Range("C12").Value2 = 0
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = 0
$ws.Range("F12").Value2 = 0
$ws.Range("G12").Value2 = 1
$ws.Range("H12").Value2 = 0
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value2 = 0
$ws.Range("K12").Value2 = 0
$ws.Range("L12").Value2 = 1

$ws.Range("A13").Value2 = "defect prediction"
$ws.Range("B13").Value2 = 1
$ws.Range("C13").Value2 = 0
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 5
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 1
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 8

$ws.Range("A14").Value2 = "image processing"
$ws.Range("B14").Value2 = 1
$ws.Range("C14").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("F14").Value2 = 0
$ws.Range("G14").Value2 = 0
$ws.Range("H14").Value2 = 0
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 0
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 1

$ws.Range("A15").Value2 = "issue close time"
$ws.Range("B15").Value2 = 0
$ws.Range("C15").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 1
$ws.Range("F15").Value2 = 0
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 1

$ws.Range("A16").Value2 = "language model"
$ws.Range("B16").Value2 = 1
$ws.Range("C16").Value2 = 0
$ws.Range("D16").Value2 = 0
$ws.Range("E16").Value2 = 9
$ws.Range("F16").Value2 = 0
$ws.Range("G16").Value2 = 2
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 0
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 12

$ws.Range("A17").Value2 = "language processing"
$ws.Range("B17").Value2 = 4
$ws.Range("C17").Value2 = 0
$ws.Range("D17").Value2 = 1
$ws.Range("E17").Value2 = 6
$ws.Range("F17").Value2 = 0
$ws.Range("G17").Value2 = 0
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 2
$ws.Range("J17").Value2 = 0
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 13

$ws.Range("A18").Value2 = "total"
$ws.Range("B18").Value2 = 18
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = 58
$ws.Range("F18").Value2 = 1
$ws.Range("G18").Value2 = 6
$ws.Range("H18").Value2 = 3
$ws.Range("I18").Value2 = 4
$ws.Range("J18").Value2 = 2
$ws.Range("K18").Value2 = 1
$ws.Range("L18").Value2 = 98

# Apply header/label style (bold, centered, bordered) to newly created cells
$styleSource = $ws.Range("K1")
$styleSource.Copy()
$ws.Range("L1").PasteSpecial(-4122)

$styleSourceA = $ws.Range("A2")
$styleSourceA.Copy()
$ws.Range("A16:A18").PasteSpecial(-4122)

$excel.CutCopyMode = 0